# Translate the customer-feedback sheet's header row and two of its data
# columns ("Service" category values, "Email" involvement-type values) to
# Spanish. ("Atmosphere" and the numeric ratings are left untouched, since
# the source workbook doesn't translate them either.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: "Date" -> "Fecha" (also renames the Excel Table's column).
$ws.Range("A1").Value = "Fecha"

# Determine the extent of the data so we touch every row of the table.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    if ($bCell.Value2 -eq "Service") {
        $bCell.Value = "Servicio"
    }

    $dCell = $ws.Cells.Item($r, 4)
    if ($dCell.Value2 -eq "Email") {
        $dCell.Value = "Correo electrónico"
    }
}
